# Update the crypto price/volume table with the latest scraped values.
# Cells that look like plain decimal numbers (e.g. "233.33") are written
# with a leading apostrophe so Excel keeps storing them as text, matching
# the existing inline-string cells in this sheet (prices with thousands
# separators such as "37.937.77" are already unambiguous text).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).Value = '37.937.77'
$ws.Cells.Item(2, 5).Value = '  +1.67%  '
$ws.Cells.Item(3, 4).Value = '2.094.02'
$ws.Cells.Item(4, 5).Value = '  -0.04%  '
$ws.Cells.Item(5, 4).Value = '''233.33'
$ws.Cells.Item(5, 5).Value = '  -0.39%  '
$ws.Cells.Item(6, 5).Value = '  +0.38%  '
$ws.Cells.Item(7, 5).Value = '  -0.03%  '
$ws.Cells.Item(8, 4).Value = '''57.69'
$ws.Cells.Item(8, 5).Value = '  +1.13%  '
$ws.Cells.Item(9, 4).Value = '''0.389'
$ws.Cells.Item(9, 5).Value = '  +1.79%  '
$ws.Cells.Item(10, 4).Value = '''0.0782'
$ws.Cells.Item(10, 5).Value = '  +2.29%  '
$ws.Cells.Item(11, 5).Value = '  +2.81%  '
$ws.Cells.Item(12, 4).Value = '2.393.48'
$ws.Cells.Item(12, 5).Value = '  +0.57%  '
$ws.Cells.Item(13, 5).Value = '  -1.43%  '
$ws.Cells.Item(14, 4).Value = '''21.23'
$ws.Cells.Item(14, 5).Value = '  +2.23%  '
$ws.Cells.Item(15, 5).Value = '  -1.50%  '
$ws.Cells.Item(16, 4).Value = '''5.26'
$ws.Cells.Item(16, 5).Value = '  +2.43%  '
$ws.Cells.Item(17, 4).Value = '2.106.33'
$ws.Cells.Item(17, 5).Value = '  +1.50%  '
$ws.Cells.Item(18, 4).Value = '37.903.08'
$ws.Cells.Item(18, 5).Value = '  +1.59%  '
$ws.Cells.Item(19, 5).Value = '  -3.00%  '
$ws.Cells.Item(20, 4).Value = '''70.92'
$ws.Cells.Item(20, 5).Value = '  +2.09%  '
$ws.Cells.Item(21, 5).Value = '  +1.17%  '
$ws.Cells.Item(22, 4).Value = '''228.75'
$ws.Cells.Item(22, 5).Value = '  +0.88%  '
$ws.Cells.Item(23, 5).Value = '  -0.06%  '
$ws.Cells.Item(24, 5).Value = '  -0.88%  '
$ws.Cells.Item(25, 5).Value = '  +0.12%  '
$ws.Cells.Item(26, 4).Value = '''170.67'
$ws.Cells.Item(26, 5).Value = '  +1.97%  '
$ws.Cells.Item(27, 5).Value = '  +11.80%  '
$ws.Cells.Item(28, 4).Value = '''8.96'
$ws.Cells.Item(28, 5).Value = '  +2.00%  '
$ws.Cells.Item(29, 5).Value = '  -0.02%  '
$ws.Cells.Item(30, 4).Value = '''19.50'
$ws.Cells.Item(30, 5).Value = '  +2.25%  '
$ws.Cells.Item(31, 5).Value = '  +1.09%  '
$ws.Cells.Item(32, 5).Value = '  +4.04%  '
$ws.Cells.Item(33, 5).Value = '  +2.10%  '
$ws.Cells.Item(34, 4).Value = '''4.59'
$ws.Cells.Item(34, 5).Value = '  +0.63%  '
$ws.Cells.Item(35, 5).Value = '  +1.37%  '
$ws.Cells.Item(36, 2).Value = 'WEMIXToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(36, 4).Value = '''1.83'
$ws.Cells.Item(36, 5).Value = '  +3.61%  '
$ws.Cells.Item(37, 2).Value = 'RenderToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(37, 4).Value = '''3.42'
$ws.Cells.Item(37, 5).Value = '  +5.95%  '
$ws.Cells.Item(38, 5).Value = '  +0.07%  '
$ws.Cells.Item(39, 4).Value = '''5.46'
$ws.Cells.Item(39, 5).Value = '  -4.00%  '
$ws.Cells.Item(40, 5).Value = '  +6.30%  '
$ws.Cells.Item(41, 4).Value = '''2.94'
$ws.Cells.Item(41, 5).Value = '  -0.63%  '
$ws.Cells.Item(42, 4).Value = '''97.41'
$ws.Cells.Item(42, 5).Value = '  +1.08%  '
$ws.Cells.Item(43, 5).Value = '  +0.80%  '
$ws.Cells.Item(44, 4).Value = '1.455.32'
$ws.Cells.Item(44, 5).Value = '  -1.38%  '
$ws.Cells.Item(45, 5).Value = '  +0.02%  '
$ws.Cells.Item(46, 5).Value = '  +3.63%  '
$ws.Cells.Item(47, 5).Value = '  +5.04%  '
$ws.Cells.Item(49, 5).Value = '  -8.32%  '
$ws.Cells.Item(50, 4).Value = '''3.03'
$ws.Cells.Item(50, 5).Value = '  +2.08%  '
$ws.Cells.Item(51, 4).Value = '2.289.18'
$ws.Cells.Item(51, 5).Value = '  +0.89%  '
